$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force a cell to store the given value as literal text, even when the
    # value looks like a number (e.g. "1", "612134", "122.5271", "$124.5271").
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
}

function Set-PlainValue($range, $value) {
    $ws.Range($range).Value = $value
}

# --- Row 2 ---
Set-PlainValue "B2" "fregrwegrew"
Set-TextValue  "C2" "1"
Set-PlainValue "D2" " "

# --- Row 3 ---
Set-PlainValue "B3" "gregre"
Set-TextValue  "C3" "2"
Set-PlainValue "D3" " "

# --- Row 4 ---
Set-PlainValue "B4" "reifire"
Set-PlainValue "C4" " "
Set-TextValue  "D4" "0"

# --- Row 5 ---
Set-PlainValue "B5" "qqq"
Set-TextValue  "C5" "0"
Set-PlainValue "D5" " "

# --- Row 6 ---
Set-PlainValue "B6" "reifire"
Set-TextValue  "C6" "1000"
Set-PlainValue "E6" "17.05.2021"

# --- Row 7 ---
Set-PlainValue "B7" "dollar to summ"
Set-PlainValue "C7" " "
Set-TextValue  "D7" "1"
Set-PlainValue "E7" "17.05.2021"

# --- Row 8 ---
Set-PlainValue "B8" "deedqwdq"
Set-TextValue  "D8" "1"
Set-PlainValue "E8" "17.05.2021"

# --- Row 9 ---
Set-TextValue  "C9" "612134"
Set-PlainValue "D9" " "
Set-PlainValue "E9" "18.05.2021"

# --- Row 10 --- (A10 turns into a real number, unlike the rest of the column)
Set-PlainValue "A10" 9
Set-PlainValue "B10" "jdijdijd"
Set-TextValue  "C10" "612134"
Set-PlainValue "D10" " "
Set-PlainValue "E10" "18.05.2021"

# --- Row 11 ---
Set-PlainValue "B11" " "
Set-TextValue  "C11" "1225271.0"
Set-TextValue  "D11" "2.0"

# --- Row 12 ---
Set-TextValue  "B12" "10000.0"
Set-TextValue  "C12" "122.5271"

# --- Row 13 (new row) ---
# A13 must end up with the same cell format (bold/border/centered, style
# index "1") as the rest of column A. A plain Font/Border assignment would
# mint a brand new style, so copy formats from A12 (which already carries
# that format) instead.
Set-PlainValue "A13" " "
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
Set-PlainValue "B13" " "
Set-TextValue  "C13" '$124.5271'
Set-PlainValue "D13" " "
Set-PlainValue "E13" " "
